# chore: adapt column header formatting to respective input file names
#
# - Rename the diff-table header suffixes from generic "_old"/"_new" to the
#   concrete format-version names "_FV2404"/"_FV2410".
# - Turn the A1:U64 range into a real Excel Table ("Table1") with autofilter.
# - Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..J ("*_old" -> "*_FV2404")
$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

# Column K ("diff") is left untouched.

# Columns L..U ("*_new" -> "*_FV2410")
$headersFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2404[$i]
}

for ($i = 0; $i -lt $headersFV2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2410[$i]
}

# Turn the used range into an Excel Table so the renamed headers become the
# table's column names, with an autofilter on the header row.
$tableRange = $ws.Range("A1:U64")
$listObject = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# Freeze the header row (split below row 1, keep row 1 visible while
# scrolling through the data rows).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
